# Add a new arrival row (row 8) to the "Main Data" sheet of BZG_Arrivals.xlsx,
# matching the extra flight record appended in the source data refresh
# ("Data downloaded from 11 airports :)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 7.0
$ws.Range("B8").Value = "Sunday, Jan 08"
$ws.Range("C8").Value = "9:20 PM"
$ws.Range("D8").Value = "FR6623"
$ws.Range("E8").Value = "London"
$ws.Range("F8").Value = "(LTN)"
$ws.Range("G8").Value = "Ryanair "
$ws.Range("H8").Value = "B738"
$ws.Range("I8").Value = "(EI-EBZ)"
$ws.Range("J8").Value = "9:07 PM"
# K column stays blank for every data row (mirrors K2:K7) - touch it so a
# (blank) cell is materialized in the sheet, just like the existing rows.
$ws.Range("K8").Borders.LineStyle = -4142
$ws.Range("L8").Value = "0 hours, -13 minutes"
# Same story for M (mirrors M2:M7).
$ws.Range("M8").Borders.LineStyle = -4142

Write-Host "Added row 8 to Main Data sheet"
